# Add "fullRNASEQ" as the library preparer's stated purpose for this run.
#
# The "purpose" column (E) currently repeats the libraryPreparer value
# ("S.GISH") for every data row; change it to "fullRNASEQ" for all rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E22").Value = "fullRNASEQ"

# Turn on iterative calculation (maximum change 1E-4), matching the
# workbook-level calculation settings saved with this edit.
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001

# Leave the selection where it ended up after editing the column, one row
# below the last data row.
[void]$ws.Range("E23").Select()
